$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "requirements" sheet, placing the copy immediately
#        before the original, producing: requirements (2), requirements, officials
$wsReq = $wb.Worksheets.Item("requirements")
$wsReq.Copy($wsReq)

$wsNew = $wb.Worksheets.Item(1)
$wsReq = $wb.Worksheets.Item("requirements")
$wsOff = $wb.Worksheets.Item("officials")

# --- 2. The native sheet-copy does not bring the table object along, so
#        drop the (vestigial, empty) table collection on the copy and
#        rebuild a 12-column table over it from scratch.
$wsNew.Range("A1").Value = "Discipline"
$wsNew.Range("B1").Value = "Vision"
$wsNew.Range("C1").Value = "Introduction to Officiating"
$wsNew.Range("D1").Value = "Rules and Regulations Part 1"
$wsNew.Range("E1").Value = "Rules and Regulations Part 2"
$wsNew.Range("F1").Value = "Rules and Regulations Part 3"
$wsNew.Range("G1").Value = "Level One Roving Umpire"
$wsNew.Range("H1").Value = "Level Two Roving Umpire"
$wsNew.Range("I1").Value = "Level One Chair Umpire"
$wsNew.Range("J1").Value = "Level One Referee"
$wsNew.Range("K1").Value = "Level Two Referee"
$wsNew.Range("L1").Value = "Level One Line Umprie"

# --- 4. This is now a blank template: clear out all the old 0/1 answers
$wsNew.Range("B2:L6").ClearContents()

# materialise empty, bordered cells for the newly added columns (I:L)
# so every row in the table actually has a cell in every column
$wsNew.Range("I2:L6").Borders.LineStyle = 1

$lo = $wsNew.ListObjects.Add(1, $wsNew.Range("A1:L6"), 0, 1)
$lo.Name = "Table14"
$lo.TableStyle = "TableStyleLight1"

# --- 5. Cosmetics that tagged along with the edit in the source workbook
$wsNew.Rows("1:1").RowHeight = 147
$wsNew.Columns("B:L").ColumnWidth = 3

$wsNew.Range("F16").Select()

$wsReq.Range("B2:H6").Select()

$wsOff.Range("K16").Select()

$wsNew.Activate()
